# v1.8 Changed reviewer verification for Admin Constrains
#
# This script reproduces (as closely as the COM-interop surface allows) the
# authoring changes described by the commit:
#  - On REVIEW-SHEET: close out the reviewer verification (J column) for the
#    three Admin Constrains review rows (SRS-ADM-001 / SRS-ADM-002), i.e.
#    change "open" -> "closed".
#  - On VERSION-HISTORY: append a new version-history row documenting the
#    change (v1.8, Omar Sherif, "Changed reviewer verification for Admin
#    Constrains", same date as the previous entry).
#  - Update the selection/active-sheet view state left behind by the edit.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # REVIEW-SHEET
$ws2 = $wb.Worksheets.Item(2)   # VERSION-HISTORY

# ---------------------------------------------------------------------
# 1. REVIEW-SHEET: close reviewer verification for the Admin Constrains
#    rows (SRS-ADM-001 / SRS-ADM-002 -> rows 15, 16, 17).
# ---------------------------------------------------------------------
$ws1.Range("J15").Value = "closed"
$ws1.Range("J16").Value = "closed"
$ws1.Range("J17").Value = "closed"

# ---------------------------------------------------------------------
# 2. VERSION-HISTORY: add the new v1.8 row, copying the formatting of an
#    existing single-line row (row 8) so the new row's cell styles match.
# ---------------------------------------------------------------------
$ws2.Range("A8:D8").Copy()
$ws2.Range("A10:D10").PasteSpecial(-4122)

$ws2.Range("A10").Value = "v1.8"
$ws2.Range("B10").Value = "Omar Sherif "
$ws2.Range("C10").Value = "Changed reviewer verification for Admin Constrains"
$ws2.Range("D10").Value = 45766

# ---------------------------------------------------------------------
# 3. View state: VERSION-HISTORY is no longer the selected/active tab,
#    REVIEW-SHEET becomes active with the updated selection.
# ---------------------------------------------------------------------
[void]$ws2.Activate()
[void]$ws2.Range("B10").Select()

[void]$ws1.Activate()
[void]$ws1.Range("J23:J24").Select()
